$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was updated from
# 2023-10-25 (45224) to 2023-11-03 (45233) for every data row (2..158).
$ws.Range("C2:C158").Value = 45233
